# Updated capital structure database
# Apply updated values to rows 2 and 3 of the earnings_debt sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    $ws.Range("D$r").Value = 0.0795
    $ws.Range("E$r").ClearContents()

    $ws.Range("G$r").Value = 0.3963585434173669
    $ws.Range("H$r").Value = 0.3963585434173669
    $ws.Range("I$r").Value = 0.1708683473389356
    $ws.Range("J$r").Value = 0.1308634564937086
    $ws.Range("K$r").Value = 9.630000000000001
    $ws.Range("L$r").Value = 0.1348739495798319
    $ws.Range("M$r").Value = 0
    $ws.Range("N$r").Value = 0
    $ws.Range("O$r").Value = 0
    $ws.Range("P$r").Value = 0
    $ws.Range("Q$r").Value = 0
    $ws.Range("R$r").Value = 0
    $ws.Range("S$r").Value = 0
    $ws.Range("T$r").ClearContents()

    $ws.Range("U$r").Value = 55.3
    $ws.Range("V$r").Value = 0.7691237830319888
    $ws.Range("W$r").Value = 0.09786585365853659
    $ws.Range("X$r").Value = 0.06313290179349144
    $ws.Range("Y$r").Value = 0.03473295186504514
    $ws.Range("Z$r").Value = 1.698624922681639
    $ws.Range("AA$r").Value = 0.2222879286684777
    $ws.Range("AB$r").Value = 0.06309072590170704
    $ws.Range("AC$r").Value = 0.1591972027667707
    $ws.Range("AD$r").Value = 0.188
    $ws.Range("AF$r").Value = 0.188
    $ws.Range("AG$r").Value = -55.11199999999999
    $ws.Range("AH$r").Value = 0.002607923648873599
    $ws.Range("AI$r").Value = 0.001755565516210966
    $ws.Range("AJ$r").Value = -3.282821062663805
    $ws.Range("AK$r").Value = -1.064184753224685
    $ws.Range("AL$r").Value = 0.008
    $ws.Range("AM$r").Value = 0.008
    $ws.Range("AN$r").Value = 0.01504
    $ws.Range("AO$r").Value = 1525
    $ws.Range("AP$r").Value = -4.40896
    $ws.Range("AQ$r").Value = 1525
}
